$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Interior.Color = 5296274
Write-Host "done"
